$d = $word.ActiveDocument

# Locate the "Thuật toán:" heading paragraph and make it bold (matching
# the formatting already used on the other headings in this document),
# both for the paragraph mark itself (pPr/rPr) and for every run in it.
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -match "Thu.t to.n:") {
        $r = $p.Range
        $r.Font.Bold = 1
        $r.Font.BoldBi = 1
    }
}
